# Apply the "Group" label renaming + view-state changes described by the diff.
#
# Content change: the experimental group labels used in the "CRA030408 "
# (bulk RNA-seq) and "OMIX012014" (processed bulk RNA-seq) sheets are
# renamed from the old ND/HFD naming scheme to the new WT/Pdcd1-/- scheme:
#   ND-IgG           -> WT+vehicle
#   ND-PD-1 mAb      -> Pdcd1-/-+vehicle
#   HFD-IgG          -> WT+PA
#   HFD-PD-1 mAb     -> Pdcd1-/-+PA
#   HFD-PD-1 mAb+Met -> Pdcd1-/-+PA+Met

$wb = $excel.ActiveWorkbook

$map = @{
    "ND-IgG"           = "WT+vehicle"
    "ND-PD-1 mAb"      = "Pdcd1-/-+vehicle"
    "HFD-IgG"          = "WT+PA"
    "HFD-PD-1 mAb"     = "Pdcd1-/-+PA"
    "HFD-PD-1 mAb+Met" = "Pdcd1-/-+PA+Met"
}

# --- Sheet "CRA030408 " (2nd sheet): column B, rows 2-16 ---
$wsBulkRaw = $wb.Worksheets.Item("CRA030408 ")
for ($r = 2; $r -le 16; $r++) {
    $cell = $wsBulkRaw.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $cell.Value = $map[$old]
    }
}

# --- Sheet "OMIX012014" (4th sheet): column B, rows 5-19 ---
$wsBulkProcessed = $wb.Worksheets.Item("OMIX012014")
for ($r = 5; $r -le 19; $r++) {
    $cell = $wsBulkProcessed.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $cell.Value = $map[$old]
    }
}

# --- View-state changes ---

# Sheet "OMIX012014": selection moves from D4 to B21
$wsBulkProcessed.Range("B21").Select()

# Sheet "CRA030408 " becomes the active/selected tab, zoom drops from 190 to
# 160, and the selection moves from B33 to D19.
$wsBulkRaw.Activate()
$excel.ActiveWindow.Zoom = 160
$wsBulkRaw.Range("D19").Select()
